# Adapt column header formatting to respective input file names:
# rename the "_old"/"_new" suffixed headers in row 1 to "_FV2304"/"_FV2310",
# wrap the sheet's data range in a table (Table1) using the new header
# names, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column headers in row 1 (A1:U1) so that the old "_old"/"_new"
#    suffixes become "_FV2304"/"_FV2310" respectively ("diff" stays as-is).
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# 2) Turn the used range into a table ("Table1") covering A1:U70, complete
#    with an autofilter on the header row, using the renamed headers above.
$dataRange = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split below row 1, keep column A as the left
#    edge) so the header stays visible while scrolling.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
